$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("O")
$ws.Columns("D").Insert(-4161, 0)
